$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 12:14"

# Rumania overtakes Belgica (rows 31/32 swap ranking position)
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 152403
$ws.Range("C31").Value = 3517
$ws.Range("D31").Value = 116628
$ws.Range("E31").Value = 30417
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 59
$ws.Range("H31").Value = 5358

$ws.Range("A32").Value = "Belgica"
$ws.Range("B32").Value = 148981
$ws.Range("C32").Value = 5385
$ws.Range("D32").Value = 20072
$ws.Range("E32").Value = 118758
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 25
$ws.Range("H32").Value = 10151

# Banglades data update (row 19)
$ws.Range("B19").Value = 377073
$ws.Range("C19").Value = 1203
$ws.Range("D19").Value = 291365
$ws.Range("E19").Value = 80208
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 5500

# El Salvador data update (row 79)
$ws.Range("E79").Value = 4069
$ws.Range("G79").Value = 6
$ws.Range("H79").Value = 887

# Finlandia overtakes Namibia (rows 101/102 swap ranking position)
$ws.Range("A101").Value = "Finlandia"
$ws.Range("B101").Value = 11849
$ws.Range("C101").Value = 269
$ws.Range("D101").Value = 8500
$ws.Range("E101").Value = 3003
$ws.Range("H101").Value = 346

$ws.Range("A102").Value = "Namibia"
$ws.Range("B102").Value = 11829
$ws.Range("D102").Value = 9778
$ws.Range("E102").Value = 1924
$ws.Range("H102").Value = 127

# Row 105 data update
$ws.Range("B105").Value = 10841
$ws.Range("C105").Value = 6
$ws.Range("E105").Value = 323

# Row 127 data update (Hong Kong)
$ws.Range("B127").Value = 5176
$ws.Range("C127").Value = 6
$ws.Range("E127").Value = 165
